$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 6557
$ws1.Range("G2").Value = 58.5
$ws1.Range("F5").Value = 411
$ws1.Range("G6").Value = 55
$ws1.Range("F14").Value = 1118
$ws1.Range("F15").Value = 3247
$ws1.Range("F17").Value = 206
$ws1.Range("F18").Value = 1895

# Sheet "全部类型" (all types) - same events, shifted by one row because of
# an extra "演出" (performance) row inserted at row 8.
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 6557
$ws4.Range("G2").Value = 58.5
$ws4.Range("F5").Value = 411
$ws4.Range("G6").Value = 55
$ws4.Range("F15").Value = 1118
$ws4.Range("F16").Value = 3247
$ws4.Range("F18").Value = 206
$ws4.Range("F19").Value = 1895
